# The author filled in the previously-blank "1024x1024 / 400x400" results
# table (rows 71-74, columns D:F, I:K, N:P) on Sheet1 with the measured
# timings that feed charts 8, 11 and 12 (Parallel For / Threads w/Mutex /
# Sequential series for the three "Samples per Pixel" groups: 9 spheres,
# 12 spheres, 15 spheres).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "9 spheres" block (columns D:F) ---
$ws.Range("D71").Value = 25766.5
$ws.Range("E71").Value = 12029.3
$ws.Range("F71").Value = 12642.7

$ws.Range("D72").Value = 108437.6
$ws.Range("E72").Value = 55705.5
$ws.Range("F72").Value = 56926.6

$ws.Range("D73").Value = 415926.2
$ws.Range("E73").Value = 217248.6
$ws.Range("F73").Value = 229301.7

$ws.Range("D74").Value = 1683655.8
$ws.Range("E74").Value = 874755.75
$ws.Range("F74").Value = 955163.4

# --- "12 spheres" block (columns I:K) ---
$ws.Range("I71").Value = 26567.9
$ws.Range("J71").Value = 13498.3
$ws.Range("K71").Value = 13601.5

$ws.Range("I72").Value = 105660.5
$ws.Range("J72").Value = 50764.1
$ws.Range("K72").Value = 53075.7

$ws.Range("I73").Value = 426819.2
$ws.Range("J73").Value = 204060.4
$ws.Range("K73").Value = 219116.79999999999

$ws.Range("I74").Value = 1652563
$ws.Range("J74").Value = 805825
$ws.Range("K74").Value = 892637.3

# --- "15 spheres" block (columns N:P) ---
$ws.Range("N71").Value = 38280.699999999997
$ws.Range("O71").Value = 15527.5
$ws.Range("P71").Value = 15916

$ws.Range("N72").Value = 136265.79999999999
$ws.Range("O72").Value = 58945
$ws.Range("P72").Value = 62556.1

$ws.Range("N73").Value = 564249.9
$ws.Range("O73").Value = 239772.3
$ws.Range("P73").Value = 258207.5

$ws.Range("N74").Value = 2159342
$ws.Range("O74").Value = 952495
$ws.Range("P74").Value = 1022660.7

# The author also scrolled/selected near the bottom of the newly filled
# table before saving.
$ws.Range("P75").Select()
